$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as literal text (not auto-converted to numbers),
# matching the source data which stores prices/changes as inline strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.104.31"
$ws.Range("E2").Value = "  +3.30%  "

$ws.Range("D3").Value = "3.061.51"
$ws.Range("E3").Value = "  +2.14%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "560.77"
$ws.Range("E5").Value = "  +3.51%  "

$ws.Range("D6").Value = "144.04"
$ws.Range("E6").Value = "  +3.95%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "3.062.16"
$ws.Range("E8").Value = "  +2.33%  "

$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  +5.09%  "

$ws.Range("E10").Value = "  +5.94%  "

$ws.Range("E11").Value = "  -9.84%  "

$ws.Range("D12").Value = "0.485"
$ws.Range("E12").Value = "  +9.47%  "

$ws.Range("E13").Value = "  +5.59%  "

$ws.Range("D14").Value = "35.48"
$ws.Range("E14").Value = "  +4.92%  "

$ws.Range("D15").Value = "3.560.32"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("D16").Value = "64.123.30"
$ws.Range("E16").Value = "  +3.20%  "

$ws.Range("D17").Value = "3.070.36"
$ws.Range("E17").Value = "  +2.46%  "

$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("E19").Value = "  +3.45%  "

$ws.Range("D20").Value = "478.50"
$ws.Range("E20").Value = "  +3.03%  "

$ws.Range("E21").Value = "  +4.95%  "

$ws.Range("D22").Value = "0.681"
$ws.Range("E22").Value = "  +4.87%  "

$ws.Range("D23").Value = "7.56"
$ws.Range("E23").Value = "  +5.46%  "

$ws.Range("D24").Value = "14.27"
$ws.Range("E24").Value = "  +14.22%  "

$ws.Range("D25").Value = "82.21"
$ws.Range("E25").Value = "  +3.91%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  +3.84%  "

$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +6.26%  "

$ws.Range("E29").Value = "  +2.18%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").Value = "26.35"
$ws.Range("E31").Value = "  +4.04%  "

$ws.Range("E32").Value = "  +1.91%  "

$ws.Range("E33").Value = "  +5.21%  "

$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  +3.58%  "

$ws.Range("D35").Value = "6.26"
$ws.Range("E35").Value = "  +7.88%  "

$ws.Range("D36").Value = "54.97"
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("E37").Value = "  +5.23%  "

$ws.Range("D38").Value = "444.52"
$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("D39").Value = "0.0812"
$ws.Range("E39").Value = "  +0.96%  "

$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +11.63%  "

$ws.Range("D41").Value = "2.992.25"
$ws.Range("E41").Value = "  +2.12%  "

$ws.Range("E42").Value = "  +2.83%  "

$ws.Range("E43").Value = "  +0.75%  "

$ws.Range("D44").Value = "27.88"
$ws.Range("E44").Value = "  +4.86%  "

$ws.Range("D45").Value = "0.262"
$ws.Range("E45").Value = "  +6.62%  "

$ws.Range("D46").Value = "2.17"
$ws.Range("E46").Value = "  +9.53%  "

$ws.Range("E48").Value = "  +4.66%  "

$ws.Range("E51").Value = "  +4.08%  "

# Row 49 <-> Row 50 swap: PEPE and Monero traded ranking positions.
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "118.64"
$ws.Range("E49").Value = "  +3.22%  "

$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0520"
$ws.Range("E50").Value = "  +5.11%  "

# Drop the transient Text number-format so the saved styles match the original (no extra style).
$ws.Range("D2:D51").ClearFormats()

Write-Output "cryptos sheet updated"
